$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1947.6
$ws.Range("I40").Value = 1784.5151
$ws.Range("K40").Value = 1784.5151
$ws.Range("M40").Value = -1609.5151

$ws.Range("H62").Value = 2749.5
$ws.Range("I62").Value = 499
$ws.Range("K62").Value = 499
$ws.Range("M62").Value = 125

$ws.Range("H65").Value = 2749.5
$ws.Range("I65").Value = 499
$ws.Range("K65").Value = 2495
$ws.Range("M65").Value = 625

$ws.Range("H74").Value = 143798.64
$ws.Range("I74").Value = 157678.5
$ws.Range("J74").Value = 5000
$ws.Range("K74").Value = 157678.5
$ws.Range("L74").Value = 5000
$ws.Range("M74").Value = -156742.5
$ws.Range("N74").Value = -6872

$ws.Range("H77").Value = 143798.64
$ws.Range("I77").Value = 157678.5
$ws.Range("J77").Value = 5000
$ws.Range("K77").Value = 788392.5
$ws.Range("L77").Value = 25000
$ws.Range("M77").Value = -783712.5
$ws.Range("N77").Value = -34360

$ws.Range("H86").Value = 1668
$ws.Range("I86").Value = 1469.125
$ws.Range("J86").Value = 1933.1666
$ws.Range("K86").Value = 1469.125
$ws.Range("L86").Value = 1933.1666
$ws.Range("M86").Value = -346.125
$ws.Range("N86").Value = -4179.1666

$ws.Range("H88").Value = 1321.15
$ws.Range("J88").Value = 1164.0714
$ws.Range("L88").Value = 1164.0714
$ws.Range("N88").Value = -1976.0714

$ws.Range("H89").Value = 1668
$ws.Range("I89").Value = 1469.125
$ws.Range("J89").Value = 1933.1666
$ws.Range("K89").Value = 7345.625
$ws.Range("L89").Value = 9665.833000000001
$ws.Range("M89").Value = -1729.625
$ws.Range("N89").Value = -20897.833

$ws.Range("H91").Value = 1321.15
$ws.Range("J91").Value = 1164.0714
$ws.Range("L91").Value = 1164.0714
$ws.Range("N91").Value = -3972.0714

$ws.Range("H100").Value = 1038.4
$ws.Range("I100").Value = 1073.75
$ws.Range("J100").Value = 897
$ws.Range("K100").Value = 1073.75
$ws.Range("L100").Value = 897
$ws.Range("M100").Value = -532.75
$ws.Range("N100").Value = -1979

$ws.Range("H125").Value = 6543.75
$ws.Range("I125").Value = 735
$ws.Range("K125").Value = 6615
$ws.Range("M125").Value = -4155

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 250
$ws.Range("I4").Value = 250
$ws.Range("K4").Value = 250
$ws.Range("M4").Value = -134

$ws.Range("H32").Value = 4533382
$ws.Range("I32").Value = 4379218.5
$ws.Range("K32").Value = 4379218.5
$ws.Range("M32").Value = -4378931.5

$ws.Range("H74").Value = 3017.5
$ws.Range("I74").Value = 2982
$ws.Range("J74").Value = 3029.3333
$ws.Range("K74").Value = 2982
$ws.Range("L74").Value = 3029.3333
$ws.Range("M74").Value = -2108
$ws.Range("N74").Value = -4777.3333

$ws.Range("H77").Value = 3017.5
$ws.Range("I77").Value = 2982
$ws.Range("J77").Value = 3029.3333
$ws.Range("K77").Value = 14910
$ws.Range("L77").Value = 15146.6665
$ws.Range("M77").Value = -10542
$ws.Range("N77").Value = -23882.6665

$ws.Range("H88").Value = 2656.3333
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 2656.3333
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 2656.3333
$ws.Range("N88").Value = -3468.3333
$ws.Range("M88").ClearContents()

$ws.Range("H91").Value = 2656.3333
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 2656.3333
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 2656.3333
$ws.Range("N91").Value = -5464.3333
$ws.Range("M91").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 485
$ws.Range("J64").Value = 179
$ws.Range("L64").Value = 179
$ws.Range("N64").Value = -629

$ws.Range("H67").Value = 485
$ws.Range("J67").Value = 179
$ws.Range("L67").Value = 179
$ws.Range("N67").Value = -1739

$ws.Range("H107").Value = 2327.75
$ws.Range("I107").Value = 3155.5
$ws.Range("K107").Value = 3155.5
$ws.Range("M107").Value = -1235.5

$ws.Range("H134").Value = 1500
$ws.Range("I134").Value = 1500
$ws.Range("K134").Value = 4500
$ws.Range("M134").Value = -1965

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3929.2
$ws.Range("J16").Value = 3499.8572
$ws.Range("L16").Value = 3499.8572
$ws.Range("N16").Value = -4073.8572

$ws.Range("H19").Value = 30
$ws.Range("I19").Value = 30
$ws.Range("K19").Value = 30
$ws.Range("M19").Value = 140

$ws.Range("H24").Value = 30
$ws.Range("I24").Value = 30
$ws.Range("K24").Value = 30
$ws.Range("M24").Value = 140

$ws.Range("H99").Value = 2069.8
$ws.Range("I99").Value = 1399
$ws.Range("J99").Value = 3076
$ws.Range("K99").Value = 1399
$ws.Range("L99").Value = 3076
$ws.Range("M99").Value = 99
$ws.Range("N99").Value = -6072

$ws.Range("H113").Value = 3929.2
$ws.Range("J113").Value = 3499.8572
$ws.Range("L113").Value = 3499.8572
$ws.Range("N113").Value = -7839.8572

$ws.Range("H126").Value = 2069.8
$ws.Range("I126").Value = 1399
$ws.Range("J126").Value = 3076
$ws.Range("K126").Value = 4197
$ws.Range("L126").Value = 9228
$ws.Range("M126").Value = -1727
$ws.Range("N126").Value = -14168

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 234
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 234
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 702
$ws.Range("N5").Value = -926
$ws.Range("M5").ClearContents()

$ws.Range("H44").Value = 419.92856
$ws.Range("I44").Value = 209.8
$ws.Range("J44").Value = 536.6667
$ws.Range("K44").Value = 629.4000000000001
$ws.Range("L44").Value = 1610.0001
$ws.Range("M44").Value = -231.4000000000001
$ws.Range("N44").Value = -2406.0001

$ws.Range("H135").Value = 234
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 234
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 2106
$ws.Range("N135").Value = -7176
$ws.Range("M135").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 1906.25
$ws.Range("J126").Value = 1750
$ws.Range("L126").Value = 5250
$ws.Range("N126").Value = -10190

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 9654.625
$ws.Range("I7").Value = 9696.799999999999
$ws.Range("K7").Value = 9696.799999999999
$ws.Range("M7").Value = -9584.799999999999

$ws.Range("H16").Value = 1569.5
$ws.Range("I16").Value = 1792.25
$ws.Range("J16").Value = 1346.75
$ws.Range("K16").Value = 1792.25
$ws.Range("L16").Value = 1346.75
$ws.Range("M16").Value = -1622.25
$ws.Range("N16").Value = -1686.75

$ws.Range("H40").Value = 5074.125
$ws.Range("I40").Value = 5013.2856
$ws.Range("J40").Value = 5500
$ws.Range("K40").Value = 5013.2856
$ws.Range("L40").Value = 5500
$ws.Range("M40").Value = -4877.2856
$ws.Range("N40").Value = -5772

$ws.Range("H122").Value = 6987.8887
$ws.Range("I122").Value = 4798.5
$ws.Range("J122").Value = 7613.4287
$ws.Range("K122").Value = 14395.5
$ws.Range("L122").Value = 22840.2861
$ws.Range("M122").Value = -11945.5
$ws.Range("N122").Value = -27740.2861

$ws.Range("H126").Value = 9654.625
$ws.Range("I126").Value = 9696.799999999999
$ws.Range("K126").Value = 29090.4
$ws.Range("M126").Value = -26620.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").ClearContents()
